# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Only column G ("K" - strikeouts) values change for data rows 2-15 on Sheet1.
# These are regenerated/computed values (no longer straight "Strike#" counts),
# so we write the recalculated literal values directly into the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 1
    6  = 1
    7  = 2
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 2
    15 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
